# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2 of the
# zh-cn and de-de worksheets to reflect the latest handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 14:50:35"
$wsZhCn.Range("H2").Value = "2016-03-12 14:50:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 14:50:38"
$wsDeDe.Range("H2").Value = "2016-03-12 14:50:58"
